$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new 20-minute trade row (row 4) with the same shape as row 3.
$ws.Range("A4").Value = 9987
$ws.Range("B4").Value = 10002
$ws.Range("C4").Value = 80.45
$ws.Range("D4").Value = 80.569999999999993
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 0.15
$ws.Range("G4").Value = 42608.624074074076
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"
$ws.Range("H4").Value = $false
